$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the run " ser necessariamente via telefone e " into two runs:
#    " ser " and "necessariamente via telefone e ".
#    A zero-length bookmark forces Word to break the run at that point; the
#    bookmark is then removed so only the run split remains.
# ---------------------------------------------------------------------------
$rngSplit = $d.Content
$rngSplit.Find.Execute(" ser necessariamente via telefone e ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $rngSplit.Start + 5
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplit", $splitRange) | Out-Null
$d.Bookmarks("TempSplit").Delete()

# ---------------------------------------------------------------------------
# 2) Fix the capitalisation of "whatsapp" -> "WhatsApp".
# ---------------------------------------------------------------------------
$rngWord = $d.Content
$rngWord.Find.Execute("whatsapp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngWord.Text = "WhatsApp"

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark (last-edit-position marker) from the empty
#    trailing paragraph to right after "WhatsA" inside "WhatsApp" - this is
#    where the user's cursor was left after typing the capital "A".
#    Bookmarks.Add with the existing "_GoBack" name relocates it.
# ---------------------------------------------------------------------------
$rngGoBack = $d.Content
$rngGoBack.Find.Execute("WhatsApp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $rngGoBack.Start + 6
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# 4) Swap the theme's dk1/lt1 colors (Text 1 <-> Background 1).
# ---------------------------------------------------------------------------
$scheme = $d.DocumentTheme.ThemeColorScheme
$dk1 = $scheme.Item(1)
$lt1 = $scheme.Item(2)
$dk1Rgb = $dk1.RGB
$lt1Rgb = $lt1.RGB
$dk1.RGB = $lt1Rgb
$lt1.RGB = $dk1Rgb
